$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("No Technique", "Naive Bayes", 1),
    @("No Technique", "Naive Bayes", 1),
    @("No Technique", "Naive Bayes", 1),
    @("No Technique", "Naive Bayes", 0.9665831244778612),
    @("No Technique", "Naive Bayes", 1),
    @("No Technique", "SVM", 0.9665831244778612),
    @("No Technique", "SVM", 0.9665831244778612),
    @("No Technique", "SVM", 1),
    @("No Technique", "SVM", 1),
    @("No Technique", "SVM", 0.9665831244778612)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
